$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 2 ---
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.1856783333333334
$ws.Cells.Item(2, 14).Value = 0.5570350000000001
$ws.Cells.Item(2, 15).Value = 0.1687333435516107
$ws.Cells.Item(2, 16).Value = 0.1687333435516107
$ws.Cells.Item(2, 17).Value = 0.020844806735
$ws.Cells.Item(2, 18).Value = 0.187603260615
$ws.Cells.Item(2, 19).Value = 0.1687333435516107
$ws.Cells.Item(2, 20).Value = 0.1687333435516107

# --- Update row 3 ---
$ws.Cells.Item(3, 15).Value = 0.04089845314263524
$ws.Cells.Item(3, 16).Value = 0.04089845314263523
$ws.Cells.Item(3, 19).Value = 0.04089845314263524
$ws.Cells.Item(3, 20).Value = 0.04089845314263523

# --- Update row 4 ---
$ws.Cells.Item(4, 13).Value = 0.8184133333333333
$ws.Cells.Item(4, 14).Value = 2.45524
$ws.Cells.Item(4, 15).Value = 0.7437249982885396
$ws.Cells.Item(4, 16).Value = 0.7437249982885394
$ws.Cells.Item(4, 17).Value = 0.09187753604
$ws.Cells.Item(4, 18).Value = 0.8268978243599999
$ws.Cells.Item(4, 19).Value = 0.7437249982885396
$ws.Cells.Item(4, 20).Value = 0.7437249982885394

# --- Add new row 5 ---
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Ucn"
$ws.Cells.Item(5, 3).Value = "Crhr2"
$ws.Cells.Item(5, 4).Value = "Resolving-Mac"
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = 0.3333333333333333
$ws.Cells.Item(5, 7).Value = 0.112263
$ws.Cells.Item(5, 8).Value = 0.336789
$ws.Cells.Item(5, 9).Value = 1
$ws.Cells.Item(5, 10).Value = 1
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.05132733333333334
$ws.Cells.Item(5, 14).Value = 0.153982
$ws.Cells.Item(5, 15).Value = 0.04664320501721457
$ws.Cells.Item(5, 16).Value = 0.04664320501721456
$ws.Cells.Item(5, 17).Value = 0.005762160422000001
$ws.Cells.Item(5, 18).Value = 0.051859443798
$ws.Cells.Item(5, 19).Value = 0.04664320501721457
$ws.Cells.Item(5, 20).Value = 0.04664320501721456
